$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'247.97"
$ws.Range("D3").Value = "'21.81"
$ws.Range("D4").Value = "'5.479"
$ws.Range("D5").Value = "'0.05644"
$ws.Range("D6").Value = "'3.380"
$ws.Range("D7").Value = "'6.436"
$ws.Range("D8").Value = "'0.8021"
$ws.Range("D9").Value = "'1.037"
$ws.Range("D11").Value = "'0.07237"
$ws.Range("D12").Value = "'0.03138"
$ws.Range("D13").Value = "'0.02939"
$ws.Range("D14").Value = "'0.09289"
$ws.Range("D15").Value = "'0.001661"
$ws.Range("D16").Value = "'3.220"
$ws.Range("D17").Value = "'0.04737"
$ws.Range("D18").Value = "'0.0005816"
$ws.Range("D19").Value = "'0.006410"
$ws.Range("D20").Value = "'0.005033"
$ws.Range("D21").Value = "'0.001050"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'0.0003203"
$ws.Range("D24").Value = "'4.164"
$ws.Range("D25").Value = "'2.108"
$ws.Range("D40").Value = "'0.04081"
$ws.Range("D41").Value = "'0.1045"
$ws.Range("D42").Value = "'0.002973"
$ws.Range("D43").Value = "'0.003267"
$ws.Range("D44").Value = "'0.009279"
$ws.Range("D45").Value = "'0.00005836"
$ws.Range("D46").Value = "'0.00000000751"
$ws.Range("D47").Value = "'0.7861"
$ws.Range("D48").Value = "'0.01663"
$ws.Range("D49").Value = "'0.00002102"
